# Refresh the cryptos price/volume table (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.539.32"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "1.587.38"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.86"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.26"
$ws.Range("E8").Value = "  +6.01%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0600"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "1.814.18"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "1.574.60"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.73"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "28.562.76"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.04"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.06"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.70"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.24"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.54"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "1.397.06"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("E36").Value = "  -10.08%  "
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E38").Value = "  +10.50%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.541"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.980"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.89"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "1.725.08"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.20"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0522"
$ws.Range("E51").Value = "  -0.82%  "
